$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2373-2414: the last two stock/index blocks (^BSESN, ^NSEMDCP50)
# that used to sit at rows 2403-2414 are moved up to the front (rows 2373-2384),
# and the remaining blocks (ASTRAL.NS, GRANULES.NS, GLENMARK.NS, BSOFT.NS, BEL.NS)
# shift down by 12 rows to occupy rows 2385-2414.
$newData = @(
    @("^BSESN", 45180, 67927.2265625, 66735.84375, 67838.6328125, 1, 67927.2265625),
    @("^BSESN", 45306, 73427.59375, 70665.5, 71683.2265625, 1, 73427.59375),
    @("^BSESN", 45390, 75124.28125, 74189.3125, 74244.8984375, 1, 75124.28125),
    @("^BSESN", 45005, 58418.78125, 57084.91015625, 57527.1015625, 2, 57084.91015625),
    @("^BSESN", 45166, 65473.26953125, 64723.62890625, 65387.16015625, 2, 64723.62890625),
    @("^BSESN", 45222, 65453.921875, 63092.98046875, 63782.80078125, 2, 63092.98046875),
    @("^NSEMDCP50", 44907, 9086.75, 8824.599609375, 8831.400390625, 1, 9086.75),
    @("^NSEMDCP50", 45180, 11789.9501953125, 11271.349609375, 11635.900390625, 1, 11789.9501953125),
    @("^NSEMDCP50", 45341, 14141, 13702.5498046875, 14062.2001953125, 1, 14141),
    @("^NSEMDCP50", 45012, 8490.900390625, 8194.599609375, 8466.7998046875, 2, 8194.599609375),
    @("^NSEMDCP50", 45222, 11373.75, 10810.75, 11033.099609375, 2, 10810.75),
    @("^NSEMDCP50", 45369, 13364.349609375, 12837, 13329.9501953125, 2, 12837),
    @("ASTRAL.NS", 45264, 2015.449951171875, 1918.25, 1928.75, 1, 2015.449951171875),
    @("ASTRAL.NS", 45355, 2143.800048828125, 2047.650024414062, 2100.25, 1, 2143.800048828125),
    @("ASTRAL.NS", 45425, 2352, 2119.199951171875, 2193.75, 1, 2352),
    @("ASTRAL.NS", 45117, 1866.949951171875, 1772.75, 1849.849975585938, 2, 1772.75),
    @("ASTRAL.NS", 45222, 1858.5, 1773, 1826.949951171875, 2, 1773),
    @("ASTRAL.NS", 45306, 1859, 1740, 1850.75, 2, 1740),
    @("GRANULES.NS", 45047, 309.7999877929688, 296.1499938964844, 296.6499938964844, 1, 309.7999877929688),
    @("GRANULES.NS", 45138, 329.2999877929688, 313.7000122070312, 318.2999877929688, 1, 329.2999877929688),
    @("GRANULES.NS", 45348, 480.4500122070312, 453.2000122070312, 463.4500122070312, 1, 480.4500122070312),
    @("GRANULES.NS", 45152, 298.3999938964844, 287.7000122070312, 290.3500061035156, 2, 287.7000122070312),
    @("GRANULES.NS", 45222, 345.6000061035156, 318.9500122070312, 326.2999877929688, 2, 318.9500122070312),
    @("GRANULES.NS", 45425, 410.2000122070312, 389.3500061035156, 406.75, 2, 389.3500061035156),
    @("GLENMARK.NS", 44788, 408, 384.5499877929688, 386.3500061035156, 1, 408),
    @("GLENMARK.NS", 44914, 451, 406.8999938964844, 429.1499938964844, 1, 451),
    @("GLENMARK.NS", 45187, 880, 775, 802.7999877929688, 1, 880),
    @("GLENMARK.NS", 45243, 783.4000244140625, 722.0999755859375, 775.9500122070312, 2, 722.0999755859375),
    @("GLENMARK.NS", 45334, 891, 771, 872.25, 2, 771),
    @("GLENMARK.NS", 45425, 1048, 985.2000122070312, 1042.449951171875, 2, 985.2000122070312),
    @("BSOFT.NS", 44900, 334.7000122070312, 306.6499938964844, 309.3500061035156, 1, 334.7000122070312),
    @("BSOFT.NS", 44949, 312.6000061035156, 287.5499877929688, 290.8999938964844, 1, 312.6000061035156),
    @("BSOFT.NS", 45327, 861.8499755859375, 817.0999755859375, 833.4000244140625, 1, 861.8499755859375),
    @("BSOFT.NS", 44956, 305.75, 250.25, 267.7999877929688, 2, 250.25),
    @("BSOFT.NS", 45012, 266.5, 250.5, 261.1499938964844, 2, 250.5),
    @("BSOFT.NS", 45425, 623.7999877929688, 564.5999755859375, 614.1500244140625, 2, 564.5999755859375),
    @("BEL.NS", 44816, 114.6500015258789, 108.5999984741211, 111, 1, 114.6500015258789),
    @("BEL.NS", 44865, 112.1999969482422, 105.1500015258789, 108.9499969482422, 1, 112.1999969482422),
    @("BEL.NS", 45180, 147.1499938964844, 133.3000030517578, 135.6999969482422, 1, 147.1499938964844),
    @("BEL.NS", 45222, 135.8999938964844, 127, 132.1999969482422, 2, 127),
    @("BEL.NS", 45334, 190.3999938964844, 171.75, 188.3000030517578, 2, 171.75),
    @("BEL.NS", 45425, 260.6000061035156, 221, 258.7999877929688, 2, 221)
)

$startRow = 2373
for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $vals = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
    $ws.Cells.Item($r, 6).Value = $vals[5]
    $ws.Cells.Item($r, 7).Value = $vals[6]
}
